# Updated cluster function to separate 20 and 16 layers data
#
# Column F ("Coax" / signal mapping) was driven by a shared formula
# (F19:F67 = F_prevBlock + 16) that assumed a single contiguous 20-layer
# cluster. The data is now split so each 16-row block (0-15, 16-31,
# 32-47, 49-64) restarts its own descending sequence, and the extra
# 65th row (E67/F67, value 64) that belonged to the old contiguous
# scheme is dropped entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each block of 16 rows holds a descending run (High down to High-15).
$blocks = @(
    @{Start = 3;  High = 15},
    @{Start = 19; High = 31},
    @{Start = 35; High = 47},
    @{Start = 51; High = 64}
)

foreach ($blk in $blocks) {
    for ($i = 0; $i -lt 16; $i++) {
        $row = $blk.Start + $i
        $val = $blk.High - $i
        $ws.Range("F$row").Value = $val
    }
}

# Row 67 no longer carries the 65th (E/F) mapping entry for this cluster.
$ws.Range("E67:F67").ClearContents()

# Match the author's final cursor position/selection in the sheet.
$ws.Range("H64").Select()
